$wb = $excel.ActiveWorkbook

$wsEdu = $wb.Worksheets.Item("Education")
$wsEdu.Range("A65").Value = "SE.COM.DURS"
$wsEdu.Range("B65").Value = "Compulsory education, duration (years)"
$wsEdu.Range("F65").Value = "Compulsory"
$wsEdu.Range("G65").Value = "Duration"

$wsEdu.Range("A66").Value = "SE.LPV.PRIM.MA"
$wsEdu.Range("B66").Value = "Learning poverty: Share of Male Children at the End-of-Primary age below minimum reading proficiency adjusted by Out-of-School Children (%)"
$wsEdu.Range("F66").Value = "Learning poverty"
$wsEdu.Range("G66").Value = "Primary education"

$wsEdu.Range("A67").Value = "SE.LPV.PRIM.FE"
$wsEdu.Range("B67").Value = "Learning poverty: Share of Female Children at the End-of-Primary age below minimum reading proficiency adjusted by Out-of-School Children (%)"
$wsEdu.Range("F67").Value = "Learning poverty"
$wsEdu.Range("G67").Value = "Primary education"

$wsEdu.Range("A68").Value = "SE.LPV.PRIM"
$wsEdu.Range("B68").Value = "Learning poverty: Share of Children at the End-of-Primary age below minimum reading proficiency adjusted by Out-of-School Children (%)"
$wsEdu.Range("F68").Value = "Learning poverty"
$wsEdu.Range("G68").Value = "Primary education"

$wsEdu.Range("A69").Value = "SE.LPV.PRIM.LD.FE"
$wsEdu.Range("B69").Value = "Female pupils below minimum reading proficiency at end of primary (%). Low GAML threshold"
$wsEdu.Range("F69").Value = "Learning poverty"
$wsEdu.Range("G69").Value = "Primary education"

$wsEdu.Range("A70").Value = "SE.LPV.PRIM.SD.FE"
$wsEdu.Range("B70").Value = "Female primary school age children out-of-school (%)"
$wsEdu.Range("F70").Value = "Learning poverty"
$wsEdu.Range("G70").Value = "Primary education"

$wsEdu.Range("A71").Value = "SE.LPV.PRIM.LD"
$wsEdu.Range("B71").Value = "Pupils below minimum reading proficiency at end of primary (%). Low GAML threshold"
$wsEdu.Range("F71").Value = "Learning poverty"
$wsEdu.Range("G71").Value = "Primary education"

$wsEdu.Range("A72").Value = "SE.LPV.PRIM.SD"
$wsEdu.Range("B72").Value = "Primary school age children out-of-school (%)"
$wsEdu.Range("F72").Value = "Learning poverty"
$wsEdu.Range("G72").Value = "Primary education"

$wsEdu.Range("A73").Value = "SE.LPV.PRIM.LD.MA"
$wsEdu.Range("B73").Value = "Male pupils below minimum reading proficiency at end of primary (%). Low GAML threshold"
$wsEdu.Range("F73").Value = "Learning poverty"
$wsEdu.Range("G73").Value = "Primary education"

$wsEdu.Range("A74").Value = "SE.LPV.PRIM.SD.MA"
$wsEdu.Range("B74").Value = "Male primary school age children out-of-school (%)"
$wsEdu.Range("F74").Value = "Learning poverty"
$wsEdu.Range("G74").Value = "Primary education"

$wsEdu.Range("A75").Value = "SE.ADT.LITR.ZS"
$wsEdu.Range("B75").Value = "Literacy rate, adult total (% of people ages 15 and above)"
$wsEdu.Range("F75").Value = "Literacy rate"
$wsEdu.Range("G75").Value = "Adult (ages 15 and above)"

$wsEdu.Range("A76").Value = "SE.ADT.LITR.MA.ZS"
$wsEdu.Range("B76").Value = "Literacy rate, adult male (% of males ages 15 and above)"
$wsEdu.Range("F76").Value = "Literacy rate"
$wsEdu.Range("G76").Value = "Adult (ages 15 and above)"

$wsEdu.Range("A77").Value = "SE.ADT.LITR.FE.ZS"
$wsEdu.Range("B77").Value = "Literacy rate, adult female (% of females ages 15 and above)"
$wsEdu.Range("F77").Value = "Literacy rate"
$wsEdu.Range("G77").Value = "Adult (ages 15 and above)"

$wsEdu.Range("A78").Value = "SE.ADT.1524.LT.FE.ZS"
$wsEdu.Range("B78").Value = "Literacy rate, youth female (% of females ages 15-24)"
$wsEdu.Range("F78").Value = "Literacy rate"
$wsEdu.Range("G78").Value = "Youth (ages 15-24)"

$wsEdu.Range("A79").Value = "SE.ADT.1524.LT.FM.ZS"
$wsEdu.Range("B79").Value = "Literacy rate, youth (ages 15-24), gender parity index (GPI)"
$wsEdu.Range("F79").Value = "Literacy rate"
$wsEdu.Range("G79").Value = "Youth (ages 15-24)"

$wsEdu.Range("A80").Value = "SE.ADT.1524.LT.MA.ZS"
$wsEdu.Range("B80").Value = "Literacy rate, youth male (% of males ages 15-24)"
$wsEdu.Range("F80").Value = "Literacy rate"
$wsEdu.Range("G80").Value = "Youth (ages 15-24)"

$wsEdu.Range("A81").Value = "SE.ADT.1524.LT.ZS"
$wsEdu.Range("B81").Value = "Literacy rate, youth total (% of people ages 15-24)"
$wsEdu.Range("F81").Value = "Literacy rate"
$wsEdu.Range("G81").Value = "Youth (ages 15-24)"

$wsFin = $wb.Worksheets.Item("Financial Sector")
$wsFin.Range("A17").Value = "FS.AST.CGOV.GD.ZS"
$wsFin.Range("B17").Value = "Claims on central government, etc. (% GDP)"
$wsFin.Range("F17").Value = "Asset"
$wsFin.Range("G17").Value = "Central government"

$wsFin.Range("A18").Value = "FS.AST.DOMS.GD.ZS"
$wsFin.Range("B18").Value = "Domestic credit provided by financial sector (% of GDP)"
$wsFin.Range("F18").Value = "Asset"
$wsFin.Range("G18").Value = "Domestic"

$wsFin.Range("A19").Value = "FS.AST.DOMO.GD.ZS"
$wsFin.Range("B19").Value = "Claims on other sectors of the domestic economy (% of GDP)"
$wsFin.Range("F19").Value = "Asset"
$wsFin.Range("G19").Value = "Other domestic sectors"

$wsFin.Range("A20").Value = "FS.AST.PRVT.GD.ZS"
$wsFin.Range("B20").Value = "Domestic credit to private sector (% of GDP)"
$wsFin.Range("F20").Value = "Asset"
$wsFin.Range("G20").Value = "Private"

$wsFin.Range("A21").Value = "FB.BNK.CAPA.ZS"
$wsFin.Range("B21").Value = "Bank capital to assets ratio (%)"
$wsFin.Range("F21").Value = "Bank"
$wsFin.Range("G21").Value = "Capital to assets"

$wsFin.Range("A22").Value = "FB.AST.NPER.ZS"
$wsFin.Range("B22").Value = "Bank nonperforming loans to total gross loans (%)"
$wsFin.Range("F22").Value = "Bank (miscellaneous)"
$wsFin.Range("G22").Value = "Nonperforming"

$wsFin.Range("A23").Value = "FD.AST.PRVT.GD.ZS"
$wsFin.Range("B23").Value = "Domestic credit to private sector by banks (% of GDP)"
$wsFin.Range("F23").Value = "Deposit money banks"
$wsFin.Range("G23").Value = "Private"

$wsFin.Range("A24").Value = "FM.AST.CGOV.ZG.M3"
$wsFin.Range("B24").Value = "Claims on central government (annual growth as % of broad money)"
$wsFin.Range("F24").Value = "Monetary Survey"
$wsFin.Range("G24").Value = "Central government"

$wsFin.Range("A25").Value = "FM.AST.DOMS.CN"
$wsFin.Range("B25").Value = "Net domestic credit (current LCU)"
$wsFin.Range("F25").Value = "Monetary Survey"
$wsFin.Range("G25").Value = "Domestic"

$wsFin.Range("A26").Value = "FM.AST.NFRG.CN"
$wsFin.Range("B26").Value = "Net foreign assets (current LCU)"
$wsFin.Range("F26").Value = "Monetary Survey"
$wsFin.Range("G26").Value = "Net foreign"

$wsFin.Range("A27").Value = "FM.AST.DOMO.ZG.M3"
$wsFin.Range("B27").Value = "Claims on other sectors of the domestic economy (annual growth as % of broad money)"
$wsFin.Range("F27").Value = "Monetary Survey"
$wsFin.Range("G27").Value = "Other domestic sectors"

$wsFin.Range("A28").Value = "FM.AST.PRVT.GD.ZS"
$wsFin.Range("B28").Value = "Monetary Sector credit to private sector (% GDP)"
$wsFin.Range("F28").Value = "Monetary Survey"
$wsFin.Range("G28").Value = "Private"

$wsFin.Range("A29").Value = "FM.AST.PRVT.ZG.M3"
$wsFin.Range("B29").Value = "Claims on private sector (annual growth as % of broad money)"
$wsFin.Range("F29").Value = "Monetary Survey"
$wsFin.Range("G29").Value = "Private"
